# Update column D (Fecha) dates for rows 126-165, and append two new rows
# (166-167) at the bottom of the data table, per the weekly refresh of the
# "Vega Monumental Concepción - Cilantro" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date (serial) values for column D, rows 126-165
$dateUpdates = @{
    126 = 44609
    127 = 44609
    128 = 44252
    129 = 44252
    130 = 44271
    131 = 44271
    132 = 44420
    133 = 44420
    134 = 44336
    135 = 44336
    136 = 44231
    137 = 44231
    138 = 44565
    139 = 44565
    140 = 44334
    141 = 44334
    142 = 44280
    143 = 44280
    144 = 44362
    145 = 44362
    146 = 44365
    147 = 44365
    148 = 44567
    149 = 44567
    150 = 44553
    151 = 44553
    152 = 44490
    153 = 44490
    154 = 44462
    155 = 44462
    156 = 44264
    157 = 44264
    158 = 44330
    159 = 44330
    160 = 44257
    161 = 44257
    162 = 44299
    163 = 44299
    164 = 44285
    165 = 44285
}

foreach ($row in $dateUpdates.Keys) {
    $ws.Cells.Item($row, 4).Value = $dateUpdates[$row]
}

# Append two new rows (166 and 167) that replicate the last "Primera"/"Segunda"
# pair (previously rows 164/165) with the newest observation date (44595).
$newRows = @(
    @{ Row = 166; A = 11; B = "Vega Monumental Concepción"; C = "Bíobío"; D = 44595; E = 8; F = 100112040; G = "Cilantro"; H = "Sin especificar"; I = "Primera";  J = 200; K = 600; L = 700; M = 650; N = "`$/atado 0,5 a 1 kilo"; O = "Región de Ñuble"; P = 650; Q = 1; R = "Hortaliza" },
    @{ Row = 167; A = 11; B = "Vega Monumental Concepción"; C = "Bíobío"; D = 44595; E = 8; F = 100112040; G = "Cilantro"; H = "Sin especificar"; I = "Segunda"; J = 100; K = 500; L = 500; M = 500; N = "`$/atado 0,5 a 1 kilo"; O = "Región de Ñuble"; P = 500; Q = 1; R = "Hortaliza" }
)

$dateNumberFormat = $ws.Cells.Item(164, 4).NumberFormat

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $r.A
    $ws.Cells.Item($row, 2).Value  = $r.B
    $ws.Cells.Item($row, 3).Value  = $r.C
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 4).NumberFormat = $dateNumberFormat
    $ws.Cells.Item($row, 5).Value  = $r.E
    $ws.Cells.Item($row, 6).Value  = $r.F
    $ws.Cells.Item($row, 7).Value  = $r.G
    $ws.Cells.Item($row, 8).Value  = $r.H
    $ws.Cells.Item($row, 9).Value  = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
}

Write-Output "Edit complete"
